# Update the role holder's name in cell A2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Rob Oudman"

# Update the active selection to A2
$ws.Range("A2").Select()
